$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Strategy & Operations")

# New "AddDataPoint(S&O)" rows, extending the table with the additional
# policy / platform / position columns referenced in the commit message.
$rows = @(
    @("AddDataPoint(S&O)", "Additional Data Point"),
    @("AddDataPoint(S&O)", "Value"),
    @("AddDataPoint(S&O)", "Position Vacant"),
    @("AddDataPoint(S&O)", "Position Missing"),
    @("AddDataPoint(S&O)", "Name of the Policy"),
    @("AddDataPoint(S&O)", "Criticality (Y/N)"),
    @("AddDataPoint(S&O)", "Availability (Y/N)"),
    @("AddDataPoint(S&O)", "Name of the Technical Platform"),
    @("AddDataPoint(S&O)", "Criticality (Y/N)"),
    @("AddDataPoint(S&O)", "Availability (Y/N)")
)

$r = 7
foreach ($row in $rows) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $r = $r + 1
}

$ws.Columns.Item(2).AutoFit() | Out-Null

# This sheet becomes the active tab/selection, replacing "D&C - Construction".
$ws.Activate()
$ws.Range("B22").Select()
